$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / percentage updates (and price fields that are not
# pure numeric strings, e.g. "41.239.54" which has two decimal points).
$plainUpdates = @{
    D2 = "41.239.54"
    E2 = "  -1.68%  "
    D3 = "2.144.20"
    E3 = "  -3.02%  "
    E4 = "  -0.13%  "
    E5 = "  -1.76%  "
    E6 = "  -3.71%  "
    E7 = "  -3.71%  "
    E8 = "  -0.05%  "
    E9 = "  -5.11%  "
    E10 = "  -7.05%  "
    E11 = "  -5.67%  "
    E12 = "  -5.91%  "
    E13 = "  -3.56%  "
    E14 = "  -5.04%  "
    D15 = "2.459.16"
    E15 = "  -3.33%  "
    E16 = "  +0.82%  "
    D17 = "2.134.49"
    E17 = "  -3.27%  "
    E18 = "  -6.43%  "
    D19 = "41.046.43"
    E19 = "  -1.86%  "
    E20 = "  -4.87%  "
    E21 = "  -5.04%  "
    E22 = "  -6.80%  "
    B23 = "InternetComputer(DFINITY)"
    C23 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    E23 = "  -8.99%  "
    B24 = "BitcoinCash"
    C24 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    E24 = "  -1.63%  "
    E25 = "  -5.45%  "
    E26 = "  +0.10%  "
    E27 = "  -7.68%  "
    E28 = "  -9.06%  "
    E29 = "  -0.79%  "
    E30 = "  -4.93%  "
    E31 = "  +1.76%  "
    E32 = "  -3.95%  "
    E33 = "  +3.79%  "
    E34 = "  -5.14%  "
    E35 = "  -9.10%  "
    E36 = "  -3.79%  "
    E37 = "  -3.65%  "
    E38 = "  -0.01%  "
    E39 = "  -2.65%  "
    E40 = "  -2.97%  "
    E41 = "  -14.85%  "
    E42 = "  -7.23%  "
    E43 = "  -10.67%  "
    E44 = "  -4.65%  "
    E45 = "  -5.29%  "
    E46 = "  -4.75%  "
    E47 = "  -6.34%  "
    E48 = "  -3.00%  "
    E49 = "  -4.93%  "
    E50 = "  -3.08%  "
    E51 = "  -8.50%  "
}
foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Price fields in column D that look like plain numbers. Assigning these
# via .Value would make Excel auto-convert them to a Number (and for very
# small values like "0.0000100" even re-format/round them), which would
# lose the original text representation. Instead, write a text formula
# that evaluates to the exact string, then freeze it to a static value via
# copy / paste-special so the stored cell keeps its original (default) style.
$numericTextUpdates = @{
    D5 = "236.49"
    D6 = "0.601"
    D7 = "69.70"
    D9 = "0.569"
    D10 = "39.04"
    D11 = "0.0895"
    D12 = "53.51"
    D13 = "0.0993"
    D14 = "6.60"
    D16 = "14.29"
    D18 = "0.777"
    D20 = "0.0000100"
    D21 = "68.81"
    D22 = "5.72"
    D23 = "9.61"
    D24 = "225.02"
    D25 = "1.93"
    D27 = "10.55"
    D28 = "3.30"
    D30 = "2.15"
    D31 = "170.18"
    D32 = "19.60"
    D33 = "31.06"
    D34 = "0.0749"
    D35 = "5.07"
    D36 = "0.120"
    D38 = "4.19"
    D39 = "0.0290"
    D40 = "2.04"
    D41 = "11.72"
    D42 = "5.22"
    D43 = "57.65"
    D44 = "0.186"
    D45 = "8.19"
    D46 = "0.0955"
    D47 = "97.17"
    D51 = "2.14"
}
foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Formula = "=" + [char]34 + $numericTextUpdates[$ref] + [char]34
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
